$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 2911
$ws.Range("E2").Value = 291
$ws.Range("F2").Value = 291
$ws.Range("G2").Value = 247
$ws.Range("H2").Value = 222
$ws.Range("I2").Value = 214
$ws.Range("J2").Value = 8
$ws.Range("K2").Value = 4912
$ws.Range("L2").Value = 2539
$ws.Range("M2").Value = 2373
$ws.Range("N2").Value = 2130
$ws.Range("O2").Value = 243
$ws.Range("P2").Value = 61
$ws.Range("Q2").Value = 237
$ws.Range("R2").Value = -285
$ws.Range("S2").Value = 23
$ws.Range("T2").Value = 318
$ws.Range("U2").Value = -82
$ws.Range("V2").Value = 2214
$ws.Range("W2").Value = 10.01
$ws.Range("X2").Value = 7.62
$ws.Range("Y2").Value = 10.49
$ws.Range("Z2").Value = 4.63
$ws.Range("AA2").Value = 107.01
$ws.Range("AB2").Value = 3433.36
$ws.Range("AC2").Value = 1846
$ws.Range("AD2").Value = 27.03
$ws.Range("AE2").Value = 18410
$ws.Range("AF2").Value = 2.71
$ws.Range("AG2").Value = 250
$ws.Range("AH2").Value = 0.5
$ws.Range("AI2").Value = 13.54
$ws.Range("AJ2").Value = 11569113

# Row 3
$ws.Range("D3").Value = 3041
$ws.Range("E3").Value = 338
$ws.Range("F3").Value = 338
$ws.Range("G3").Value = 276
$ws.Range("H3").Value = 252
$ws.Range("I3").Value = 228
$ws.Range("J3").Value = 24
$ws.Range("K3").Value = 5568
$ws.Range("L3").Value = 2937
$ws.Range("M3").Value = 2630
$ws.Range("N3").Value = 2365
$ws.Range("O3").Value = 265
$ws.Range("P3").Value = 61
$ws.Range("Q3").Value = 262
$ws.Range("R3").Value = -484
$ws.Range("S3").Value = 246
$ws.Range("T3").Value = 368
$ws.Range("U3").Value = -106
$ws.Range("V3").Value = 2526
$ws.Range("W3").Value = 11.13
$ws.Range("X3").Value = 8.289999999999999
$ws.Range("Y3").Value = 10.16
$ws.Range("Z3").Value = 4.81
$ws.Range("AA3").Value = 111.66
$ws.Range("AB3").Value = 3755.68
$ws.Range("AC3").Value = 1974
$ws.Range("AD3").Value = 21.5
$ws.Range("AE3").Value = 20444
$ws.Range("AF3").Value = 2.08
$ws.Range("AG3").Value = 250
$ws.Range("AH3").Value = 0.59
$ws.Range("AI3").Value = 12.66
$ws.Range("AJ3").Value = 11569113

# Row 4
$ws.Range("D4").Value = 3404
$ws.Range("E4").Value = 401
$ws.Range("F4").Value = 401
$ws.Range("G4").Value = 337
$ws.Range("H4").Value = 307
$ws.Range("I4").Value = 273
$ws.Range("J4").Value = 34
$ws.Range("K4").Value = 5859
$ws.Range("L4").Value = 2975
$ws.Range("M4").Value = 2884
$ws.Range("N4").Value = 2589
$ws.Range("O4").Value = 295
$ws.Range("P4").Value = 61
$ws.Range("Q4").Value = 530
$ws.Range("R4").Value = -431
$ws.Range("S4").Value = -68
$ws.Range("T4").Value = 366
$ws.Range("U4").Value = 164
$ws.Range("V4").Value = 2512
$ws.Range("W4").Value = 11.79
$ws.Range("X4").Value = 9.02
$ws.Range("Y4").Value = 11.02
$ws.Range("Z4").Value = 5.37
$ws.Range("AA4").Value = 103.13
$ws.Range("AB4").Value = 4128.58
$ws.Range("AC4").Value = 2360
$ws.Range("AD4").Value = 15.74
$ws.Range("AE4").Value = 22378
$ws.Range("AF4").Value = 1.66
$ws.Range("AG4").Value = 330
$ws.Range("AH4").Value = 0.89
$ws.Range("AI4").Value = 13.98
$ws.Range("AJ4").Value = 11569113

# Row 5
$ws.Range("D5").Value = 3500
$ws.Range("E5").Value = 373
$ws.Range("F5").Value = 373
$ws.Range("G5").Value = 310
$ws.Range("H5").Value = 264
$ws.Range("I5").Value = 231
$ws.Range("J5").Value = 33
$ws.Range("K5").Value = 5953
$ws.Range("L5").Value = 2944
$ws.Range("M5").Value = 3008
$ws.Range("N5").Value = 2699
$ws.Range("O5").Value = 309
$ws.Range("P5").Value = 61
$ws.Range("Q5").Value = 457
$ws.Range("R5").Value = -435
$ws.Range("S5").Value = -72
$ws.Range("T5").Value = 466
$ws.Range("U5").Value = -9
$ws.Range("V5").Value = 2456
$ws.Range("W5").Value = 10.65
$ws.Range("X5").Value = 7.54
$ws.Range("Y5").Value = 8.73
$ws.Range("Z5").Value = 4.47
$ws.Range("AA5").Value = 97.88
$ws.Range("AB5").Value = 4408.25
$ws.Range("AC5").Value = 1995
$ws.Range("AD5").Value = 16.39
$ws.Range("AE5").Value = 23637
$ws.Range("AF5").Value = 1.38
$ws.Range("AG5").Value = 300
$ws.Range("AH5").Value = 0.92
$ws.Range("AI5").Value = 14.86
$ws.Range("AJ5").Value = 11569113

# Row 6
$ws.Range("D6").Value = 3935
$ws.Range("E6").Value = 393
$ws.Range("F6").Value = 393
$ws.Range("G6").Value = 351
$ws.Range("H6").Value = 294
$ws.Range("I6").Value = 257
$ws.Range("K6").Value = 6391
$ws.Range("L6").Value = 3309
$ws.Range("M6").Value = 3082
$ws.Range("N6").Value = 2736
$ws.Range("P6").Value = 61
$ws.Range("Q6").Value = 490
$ws.Range("R6").Value = -508
$ws.Range("S6").Value = 32
$ws.Range("T6").Value = 607
$ws.Range("U6").Value = -117
$ws.Range("V6").Value = 2734
$ws.Range("W6").Value = 9.99
$ws.Range("X6").Value = 7.47
$ws.Range("Y6").Value = 9.470000000000001
$ws.Range("Z6").Value = 4.76
$ws.Range("AA6").Value = 107.38
$ws.Range("AB6").Value = 4785.2
$ws.Range("AC6").Value = 2225
$ws.Range("AD6").Value = 10.52
$ws.Range("AE6").Value = 24607
$ws.Range("AF6").Value = 0.95
$ws.Range("AG6").Value = 350
$ws.Range("AH6").Value = 1.5
$ws.Range("AI6").Value = 15.12
$ws.Range("AJ6").Value = 11569113

# Row 7
$ws.Range("D7").Value = 4572
$ws.Range("E7").Value = 492
$ws.Range("G7").Value = 498
$ws.Range("H7").Value = 424
$ws.Range("I7").Value = 366
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 3542
$ws.Range("M7").Value = 3458
$ws.Range("N7").Value = 3102
$ws.Range("P7").Value = 61
$ws.Range("Q7").Value = 525
$ws.Range("R7").Value = -644
$ws.Range("S7").Value = 166
$ws.Range("T7").Value = 530
$ws.Range("W7").Value = 10.77
$ws.Range("X7").Value = 9.27
$ws.Range("Y7").Value = 12.54
$ws.Range("Z7").Value = 6.33
$ws.Range("AA7").Value = 102.41
$ws.Range("AC7").Value = 3164
$ws.Range("AD7").Value = 11.38
$ws.Range("AE7").Value = 27898
$ws.Range("AF7").Value = 1.29
$ws.Range("AG7").Value = 400
$ws.Range("AH7").Value = 1.11
$ws.Range("AI7").Value = 12.64

# Row 8
$ws.Range("D8").Value = 5058
$ws.Range("E8").Value = 549
$ws.Range("G8").Value = 520
$ws.Range("H8").Value = 426
$ws.Range("I8").Value = 365
$ws.Range("K8").Value = 7338
$ws.Range("L8").Value = 3560
$ws.Range("M8").Value = 3778
$ws.Range("N8").Value = 3422
$ws.Range("P8").Value = 61
$ws.Range("Q8").Value = 630
$ws.Range("R8").Value = -360
$ws.Range("S8").Value = -53
$ws.Range("T8").Value = 360
$ws.Range("W8").Value = 10.85
$ws.Range("X8").Value = 8.42
$ws.Range("Y8").Value = 11.19
$ws.Range("Z8").Value = 5.94
$ws.Range("AA8").Value = 94.2
$ws.Range("AC8").Value = 3155
$ws.Range("AD8").Value = 11.41
$ws.Range("AE8").Value = 30776
$ws.Range("AF8").Value = 1.17
$ws.Range("AG8").Value = 400
$ws.Range("AH8").Value = 1.11
$ws.Range("AI8").Value = 12.68

# Row 9
$ws.Range("D9").Value = 5604
$ws.Range("E9").Value = 614
$ws.Range("G9").Value = 585
$ws.Range("H9").Value = 476
$ws.Range("I9").Value = 412
$ws.Range("K9").Value = 7728
$ws.Range("L9").Value = 3582
$ws.Range("M9").Value = 4146
$ws.Range("N9").Value = 3789
$ws.Range("P9").Value = 61
$ws.Range("Q9").Value = 655
$ws.Range("R9").Value = -293
$ws.Range("S9").Value = -53
$ws.Range("T9").Value = 292
$ws.Range("W9").Value = 10.96
$ws.Range("X9").Value = 8.5
$ws.Range("Y9").Value = 11.44
$ws.Range("Z9").Value = 6.33
$ws.Range("AA9").Value = 86.41
$ws.Range("AC9").Value = 3566
$ws.Range("AD9").Value = 10.1
$ws.Range("AE9").Value = 34076
$ws.Range("AF9").Value = 1.06
$ws.Range("AG9").Value = 400
$ws.Range("AH9").Value = 1.11
$ws.Range("AI9").Value = 11.22

# Remove cells that are deleted entirely in the target (ClearContents removes the <c> element)
$ws.Range("U7").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("U9").ClearContents()
